# Scheduled market-data refresh: update pulled Universalis price figures
# (currentAveragePrice / NQ / HQ / Leve sell price / profit columns, H:N)
# for the rows whose prices moved since the last run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value2 = 3264.7727
$ws.Cells.Item(17, 10).Value2 = 3264.7727
$ws.Cells.Item(17, 12).Value2 = 9794.3181
$ws.Cells.Item(17, 14).Value2 = -10130.3181
$ws.Cells.Item(19, 8).Value2 = 1307.7727
$ws.Cells.Item(19, 10).Value2 = 961.875
$ws.Cells.Item(19, 12).Value2 = 961.875
$ws.Cells.Item(19, 14).Value2 = -1311.875
$ws.Cells.Item(76, 8).Value2 = 6271.75
$ws.Cells.Item(76, 9).Value2 = 5254.5454
$ws.Cells.Item(76, 10).Value2 = 7132.4614
$ws.Cells.Item(76, 11).Value2 = 5254.5454
$ws.Cells.Item(76, 12).Value2 = 7132.4614
$ws.Cells.Item(76, 13).Value2 = -4939.5454
$ws.Cells.Item(76, 14).Value2 = -7762.4614
$ws.Cells.Item(79, 8).Value2 = 6271.75
$ws.Cells.Item(79, 9).Value2 = 5254.5454
$ws.Cells.Item(79, 10).Value2 = 7132.4614
$ws.Cells.Item(79, 11).Value2 = 5254.5454
$ws.Cells.Item(79, 12).Value2 = 7132.4614
$ws.Cells.Item(79, 13).Value2 = -4162.5454
$ws.Cells.Item(79, 14).Value2 = -9316.4614
$ws.Cells.Item(80, 8).Value2 = 1856.3334
$ws.Cells.Item(80, 9).Value2 = 2469.2
$ws.Cells.Item(80, 10).Value2 = 1549.9
$ws.Cells.Item(80, 11).Value2 = 7407.599999999999
$ws.Cells.Item(80, 12).Value2 = 4649.700000000001
$ws.Cells.Item(80, 13).Value2 = -6409.599999999999
$ws.Cells.Item(80, 14).Value2 = -6645.700000000001
$ws.Cells.Item(83, 8).Value2 = 1856.3334
$ws.Cells.Item(83, 9).Value2 = 2469.2
$ws.Cells.Item(83, 10).Value2 = 1549.9
$ws.Cells.Item(83, 11).Value2 = 22222.8
$ws.Cells.Item(83, 12).Value2 = 13949.1
$ws.Cells.Item(83, 13).Value2 = -17230.8
$ws.Cells.Item(83, 14).Value2 = -23933.1
$ws.Cells.Item(132, 8).Value2 = 4620.8
$ws.Cells.Item(132, 9).Value2 = 4355
$ws.Cells.Item(132, 10).Value2 = 11000
$ws.Cells.Item(132, 11).Value2 = 13065
$ws.Cells.Item(132, 12).Value2 = 33000
$ws.Cells.Item(132, 13).Value2 = -10535
$ws.Cells.Item(132, 14).Value2 = -38060
$ws.Cells.Item(137, 8).Value2 = 3774428.2
$ws.Cells.Item(137, 9).Value2 = 811.8095
$ws.Cells.Item(137, 10).Value2 = 18182782
$ws.Cells.Item(137, 11).Value2 = 2435.4285
$ws.Cells.Item(137, 12).Value2 = 54548346
$ws.Cells.Item(137, 13).Value2 = 114.5715
$ws.Cells.Item(137, 14).Value2 = -54553446
$ws.Cells.Item(138, 8).Value2 = 2115.0645
$ws.Cells.Item(138, 10).Value2 = 2935.9062
$ws.Cells.Item(138, 12).Value2 = 8807.7186
$ws.Cells.Item(138, 14).Value2 = -19087.7186
$ws.Cells.Item(141, 8).Value2 = 2720
$ws.Cells.Item(141, 9).Value2 = 2068.5715
$ws.Cells.Item(141, 10).Value2 = 5000
$ws.Cells.Item(141, 11).Value2 = 6205.7145
$ws.Cells.Item(141, 12).Value2 = 15000
$ws.Cells.Item(141, 13).Value2 = -1025.7145
$ws.Cells.Item(141, 14).Value2 = -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value2 = 1670.9445
$ws.Cells.Item(45, 9).Value2 = 1728.6154
$ws.Cells.Item(45, 10).Value2 = 1521
$ws.Cells.Item(45, 11).Value2 = 1728.6154
$ws.Cells.Item(45, 12).Value2 = 1521
$ws.Cells.Item(45, 13).Value2 = -1351.6154
$ws.Cells.Item(45, 14).Value2 = -2275
$ws.Cells.Item(61, 8).Value2 = 21740806
$ws.Cells.Item(61, 9).Value2 = 29413610
$ws.Cells.Item(61, 11).Value2 = 29413610
$ws.Cells.Item(61, 13).Value2 = -29413398
$ws.Cells.Item(74, 8).Value2 = 15628337
$ws.Cells.Item(74, 9).Value2 = 27779746
$ws.Cells.Item(74, 11).Value2 = 27779746
$ws.Cells.Item(74, 13).Value2 = -27778872
$ws.Cells.Item(77, 8).Value2 = 15628337
$ws.Cells.Item(77, 9).Value2 = 27779746
$ws.Cells.Item(77, 11).Value2 = 138898730
$ws.Cells.Item(77, 13).Value2 = -138894362
$ws.Cells.Item(102, 8).Value2 = 2225
$ws.Cells.Item(102, 9).Value2 = 1000
$ws.Cells.Item(102, 10).Value2 = 2633.3333
$ws.Cells.Item(102, 11).Value2 = 1000
$ws.Cells.Item(102, 12).Value2 = 2633.3333
$ws.Cells.Item(102, 13).Value2 = 622
$ws.Cells.Item(102, 14).Value2 = -5877.3333
$ws.Cells.Item(122, 8).Value2 = 6025.591
$ws.Cells.Item(122, 9).Value2 = 6958.278
$ws.Cells.Item(122, 10).Value2 = 1828.5
$ws.Cells.Item(122, 11).Value2 = 20874.834
$ws.Cells.Item(122, 12).Value2 = 5485.5
$ws.Cells.Item(122, 13).Value2 = -18424.834
$ws.Cells.Item(122, 14).Value2 = -10385.5
$ws.Cells.Item(132, 8).Value2 = 11908290
$ws.Cells.Item(132, 9).Value2 = 22730532
$ws.Cells.Item(132, 10).Value2 = 3824
$ws.Cells.Item(132, 11).Value2 = 68191596
$ws.Cells.Item(132, 12).Value2 = 11472
$ws.Cells.Item(132, 13).Value2 = -68189066
$ws.Cells.Item(132, 14).Value2 = -16532
$ws.Cells.Item(136, 8).Value2 = 21740806
$ws.Cells.Item(136, 9).Value2 = 29413610
$ws.Cells.Item(136, 11).Value2 = 88240830
$ws.Cells.Item(136, 13).Value2 = -88238280

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value2 = 16668615
$ws.Cells.Item(86, 9).Value2 = 1617.2858
$ws.Cells.Item(86, 10).Value2 = 55558276
$ws.Cells.Item(86, 11).Value2 = 1617.2858
$ws.Cells.Item(86, 12).Value2 = 55558276
$ws.Cells.Item(86, 13).Value2 = -494.2858000000001
$ws.Cells.Item(86, 14).Value2 = -55560522
$ws.Cells.Item(89, 8).Value2 = 16668615
$ws.Cells.Item(89, 9).Value2 = 1617.2858
$ws.Cells.Item(89, 10).Value2 = 55558276
$ws.Cells.Item(89, 11).Value2 = 8086.429
$ws.Cells.Item(89, 12).Value2 = 277791380
$ws.Cells.Item(89, 13).Value2 = -2470.429
$ws.Cells.Item(89, 14).Value2 = -277802612
$ws.Cells.Item(99, 8).Value2 = 972.5
$ws.Cells.Item(99, 9).Value2 = 851.6667
$ws.Cells.Item(99, 11).Value2 = 851.6667
$ws.Cells.Item(99, 13).Value2 = 646.3333
$ws.Cells.Item(105, 8).Value2 = 4691.7393
$ws.Cells.Item(105, 9).Value2 = 2727.5
$ws.Cells.Item(105, 11).Value2 = 2727.5
$ws.Cells.Item(105, 13).Value2 = -980.5
$ws.Cells.Item(134, 8).Value2 = 3403.7576
$ws.Cells.Item(134, 9).Value2 = 2158.5217
$ws.Cells.Item(134, 10).Value2 = 6267.8
$ws.Cells.Item(134, 11).Value2 = 6475.5651
$ws.Cells.Item(134, 12).Value2 = 18803.4
$ws.Cells.Item(134, 13).Value2 = -3940.5651
$ws.Cells.Item(134, 14).Value2 = -23873.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 7409037
$ws.Cells.Item(31, 9).Value2 = 1666.5227
$ws.Cells.Item(31, 11).Value2 = 1666.5227
$ws.Cells.Item(31, 13).Value2 = -1371.5227
$ws.Cells.Item(34, 8).Value2 = 7409037
$ws.Cells.Item(34, 9).Value2 = 1666.5227
$ws.Cells.Item(34, 11).Value2 = 1666.5227
$ws.Cells.Item(34, 13).Value2 = -1464.5227
$ws.Cells.Item(58, 8).Value2 = 2208.1538
$ws.Cells.Item(58, 9).Value2 = 892.6667
$ws.Cells.Item(58, 10).Value2 = 3335.7144
$ws.Cells.Item(58, 11).Value2 = 892.6667
$ws.Cells.Item(58, 12).Value2 = 3335.7144
$ws.Cells.Item(58, 13).Value2 = -689.6667
$ws.Cells.Item(58, 14).Value2 = -3741.7144
$ws.Cells.Item(94, 8).Value2 = 3177.2693
$ws.Cells.Item(94, 9).Value2 = 1954.7273
$ws.Cells.Item(94, 10).Value2 = 4073.8
$ws.Cells.Item(94, 11).Value2 = 1954.7273
$ws.Cells.Item(94, 12).Value2 = 4073.8
$ws.Cells.Item(94, 13).Value2 = -1503.7273
$ws.Cells.Item(94, 14).Value2 = -4975.8
$ws.Cells.Item(122, 8).Value2 = 2154.1333
$ws.Cells.Item(122, 9).Value2 = 2226
$ws.Cells.Item(122, 10).Value2 = 1866.6666
$ws.Cells.Item(122, 11).Value2 = 6678
$ws.Cells.Item(122, 12).Value2 = 5599.9998
$ws.Cells.Item(122, 13).Value2 = -4228
$ws.Cells.Item(122, 14).Value2 = -10499.9998
$ws.Cells.Item(132, 8).Value2 = 15627016
$ws.Cells.Item(132, 9).Value2 = 17242846
$ws.Cells.Item(132, 10).Value2 = 7337.3335
$ws.Cells.Item(132, 11).Value2 = 51728538
$ws.Cells.Item(132, 12).Value2 = 22012.0005
$ws.Cells.Item(132, 13).Value2 = -51726008
$ws.Cells.Item(132, 14).Value2 = -27072.0005
$ws.Cells.Item(134, 8).Value2 = 681470.2
$ws.Cells.Item(134, 9).Value2 = 1310.4333
$ws.Cells.Item(134, 10).Value2 = 4762429
$ws.Cells.Item(134, 11).Value2 = 3931.2999
$ws.Cells.Item(134, 12).Value2 = 14287287
$ws.Cells.Item(134, 13).Value2 = -1396.2999
$ws.Cells.Item(134, 14).Value2 = -14292357
$ws.Cells.Item(136, 8).Value2 = 2208.1538
$ws.Cells.Item(136, 9).Value2 = 892.6667
$ws.Cells.Item(136, 10).Value2 = 3335.7144
$ws.Cells.Item(136, 11).Value2 = 2678.0001
$ws.Cells.Item(136, 12).Value2 = 10007.1432
$ws.Cells.Item(136, 13).Value2 = -128.0001000000002
$ws.Cells.Item(136, 14).Value2 = -15107.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value2 = 39.346153
$ws.Cells.Item(12, 9).Value2 = 29.5
$ws.Cells.Item(12, 11).Value2 = 88.5
$ws.Cells.Item(12, 13).Value2 = 84.5
$ws.Cells.Item(113, 8).Value2 = 810.6
$ws.Cells.Item(113, 9).Value2 = 487
$ws.Cells.Item(113, 10).Value2 = 972.4
$ws.Cells.Item(113, 11).Value2 = 1461
$ws.Cells.Item(113, 12).Value2 = 2917.2
$ws.Cells.Item(113, 13).Value2 = 709
$ws.Cells.Item(113, 14).Value2 = -7257.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value2 = 2201
$ws.Cells.Item(97, 9).Value2 = 2201
$ws.Cells.Item(97, 11).Value2 = 2201
$ws.Cells.Item(97, 13).Value2 = -1705
$ws.Cells.Item(102, 8).Value2 = 2469.0625
$ws.Cells.Item(102, 9).Value2 = 2670
$ws.Cells.Item(102, 11).Value2 = 2670
$ws.Cells.Item(102, 13).Value2 = -1048
$ws.Cells.Item(132, 8).Value2 = 4541.8
$ws.Cells.Item(132, 9).Value2 = 3345.45
$ws.Cells.Item(132, 10).Value2 = 6136.933
$ws.Cells.Item(132, 11).Value2 = 10036.35
$ws.Cells.Item(132, 12).Value2 = 18410.799
$ws.Cells.Item(132, 13).Value2 = -7506.349999999999
$ws.Cells.Item(132, 14).Value2 = -23470.799

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value2 = 13333967
$ws.Cells.Item(2, 10).Value2 = 13333967
$ws.Cells.Item(2, 12).Value2 = 13333967
$ws.Cells.Item(2, 14).Value2 = -13334191
$ws.Cells.Item(122, 8).Value2 = 8224.666999999999
$ws.Cells.Item(122, 9).Value2 = 11398.333
$ws.Cells.Item(122, 10).Value2 = 6108.8887
$ws.Cells.Item(122, 11).Value2 = 34194.999
$ws.Cells.Item(122, 12).Value2 = 18326.6661
$ws.Cells.Item(122, 13).Value2 = -31744.999
$ws.Cells.Item(122, 14).Value2 = -23226.6661
$ws.Cells.Item(132, 8).Value2 = 9440818
$ws.Cells.Item(132, 9).Value2 = 4139.6665
$ws.Cells.Item(132, 11).Value2 = 12418.9995
$ws.Cells.Item(132, 13).Value2 = -9888.999500000002
$ws.Cells.Item(136, 8).Value2 = 13163098
$ws.Cells.Item(136, 9).Value2 = 18520434
$ws.Cells.Item(136, 10).Value2 = 13273.182
$ws.Cells.Item(136, 11).Value2 = 55561302
$ws.Cells.Item(136, 12).Value2 = 39819.546
$ws.Cells.Item(136, 13).Value2 = -55558752

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value2 = 1712.841
$ws.Cells.Item(132, 9).Value2 = 1207.4572
$ws.Cells.Item(132, 10).Value2 = 3678.2222
$ws.Cells.Item(132, 11).Value2 = 3622.3716
$ws.Cells.Item(132, 12).Value2 = 11034.6666
$ws.Cells.Item(132, 13).Value2 = -1092.3716
$ws.Cells.Item(132, 14).Value2 = -16094.6666
$ws.Cells.Item(136, 8).Value2 = 1161
$ws.Cells.Item(136, 9).Value2 = 1234.9
$ws.Cells.Item(136, 10).Value2 = 668.3333
$ws.Cells.Item(136, 11).Value2 = 3704.7
$ws.Cells.Item(136, 12).Value2 = 2004.9999
$ws.Cells.Item(136, 13).Value2 = -1154.7
$ws.Cells.Item(136, 14).Value2 = -7104.9999
